$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (also updates the _FilterDatabase defined-name reference)
$ws.Name = "C_15.2"

# Make sure the defined name reference is quoted like Excel would emit it
# for a sheet name containing a period.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "_FilterDatabase" -or $n.Name -like "*_FilterDatabase") {
        $n.RefersTo = "='C_15.2'!`$B`$5:`$B`$5"
    }
}

# Re-apply the custom number format to the "white" (unshaded) striped data
# rows of the table so their style entry is regenerated.
$dataRng = $ws.Range("C7:H7,C9:H9,C11:H11,C13:H13,C15:H15,C17:H17,C19:H19,C21:H21")
$dataRng.NumberFormat = "#,##0.0"

